# update tien dam cuoi
# Adds a handful of new contribution entries across several sheets, moves
# three entries from "Đồng nghiệp" over to "Đại học", and leaves the
# final selection/active sheet on "Đồng nghiệp" to match the source edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Đồng nghiệp" (5th tab, sheet5.xml)
# New shared strings must be introduced here first so they end up with
# indices 44, 45, 46 like in the target workbook.
# ---------------------------------------------------------------------
$wsDongNghiep = $wb.Worksheets.Item(5)

# These three people move out to "Đại học" - clear them here first.
$wsDongNghiep.Range("A9:B11").ClearContents()

# New row 5: Đức(FPT) - 500000
$wsDongNghiep.Range("A5").Value = "Đức(FPT)"
$wsDongNghiep.Range("B5").Value = 500000

# Move the SUM total row from 13 down to 17 to make room for new entries.
$wsDongNghiep.Range("A17").Value = "SUM"
$wsDongNghiep.Range("B17").Formula = "=SUM(B1:B16)"
$wsDongNghiep.Range("A13:B13").ClearContents()

# New row 13: Tuấn Anh (MOnkey) - 300000
$wsDongNghiep.Range("A13").Value = "Tuấn Anh (MOnkey)"
$wsDongNghiep.Range("B13").Value = 300000

# New row 14: Thủy(Monkey) - 300000
$wsDongNghiep.Range("A14").Value = "Thủy(Monkey)"
$wsDongNghiep.Range("B14").Value = 300000

# ---------------------------------------------------------------------
# Sheet "Cấp3" (3rd tab, sheet3.xml) - two new rows
# ---------------------------------------------------------------------
$wsCap3 = $wb.Worksheets.Item(3)

$wsCap3.Range("A5").Value = "Trường Nhung"
$wsCap3.Range("B5").Value = 500000

$wsCap3.Range("A6").Value = "Lịch Tời"
$wsCap3.Range("B6").Value = 500000

# ---------------------------------------------------------------------
# Sheet "Anh em họ hàng" (1st tab, sheet1.xml) - two new rows
# ---------------------------------------------------------------------
$wsAnhEm = $wb.Worksheets.Item(1)

$wsAnhEm.Range("A4").Value = "Anh Cường(Bác Bộc)"
$wsAnhEm.Range("B4").Value = 500000

$wsAnhEm.Range("A5").Value = "Anh Cao()"
$wsAnhEm.Range("B5").Value = 1000000

$wsAnhEm.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------------
# Sheet "Đại học" (4th tab, sheet4.xml) - entries that moved in from
# "Đồng nghiệp" (reuses existing shared strings, no new ones created).
# ---------------------------------------------------------------------
$wsDaiHoc = $wb.Worksheets.Item(4)

$wsDaiHoc.Range("A8").Value = "Anh Toàn"
$wsDaiHoc.Range("B8").Value = 500000

$wsDaiHoc.Range("A9").Value = "Anh Sơn"

$wsDaiHoc.Range("A10").Value = "Anh Xuân"

# ---------------------------------------------------------------------
# Selections per sheet, matching the final cursor position left on each
# tab, and the workbook's active sheet/tab.
# ---------------------------------------------------------------------
$wsAnhEm.Range("A6").Select()

$wb.Worksheets.Item(2).Range("B41").Select()

$wsCap3.Range("D22").Select()

$wsDaiHoc.Range("B7").Select()

$wsDongNghiep.Activate()
$wsDongNghiep.Range("B9").Select()
